$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the category code/name columns left (kategori_id column is dropped),
# and clear out the now-unused third column.
$ws.Range("A1").Value = "kategori_kode"
$ws.Range("B1").Value = "kategori_nama"
$ws.Range("C1").ClearContents()

$ws.Range("A2").Value = "FOD"
$ws.Range("B2").Value = "Dessert"
$ws.Range("C2").ClearContents()

$ws.Range("A3").Value = "MUW"
$ws.Range("B3").Value = "Make Up"
$ws.Range("C3").ClearContents()

$ws.Range("A4").Value = "SS"
$ws.Range("B4").Value = "Skincare"
$ws.Range("C4").ClearContents()

# Match the resulting selection left behind in the saved file.
$ws.Range("C1:C4").Select()
